$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2126.1384
$ws.Range("J17").Value = 2126.1384
$ws.Range("L17").Value = 6378.415199999999
$ws.Range("N17").Value = -6714.415199999999
$ws.Range("H112").Value = 1902.8206
$ws.Range("J112").Value = 1902.8206
$ws.Range("L112").Value = 5708.4618
$ws.Range("N112").Value = -7924.4618
$ws.Range("H129").Value = 827.3099999999999
$ws.Range("J129").Value = 916.5625
$ws.Range("L129").Value = 2749.6875
$ws.Range("N129").Value = -12749.6875
$ws.Range("H137").Value = 1675660.9
$ws.Range("I137").Value = 4525730.5
$ws.Range("J137").Value = 4930.483
$ws.Range("K137").Value = 13577191.5
$ws.Range("L137").Value = 14791.449
$ws.Range("M137").Value = -13574641.5
$ws.Range("N137").Value = -19891.449
$ws.Range("H138").Value = 1410.51
$ws.Range("I138").Value = 820.34375
$ws.Range("J138").Value = 1688.2354
$ws.Range("K138").Value = 2461.03125
$ws.Range("L138").Value = 5064.706200000001
$ws.Range("M138").Value = 2678.96875
$ws.Range("N138").Value = -15344.7062
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1745.1143
$ws.Range("I2").Value = 1797.2307
$ws.Range("J2").Value = 1594.5555
$ws.Range("K2").Value = 1797.2307
$ws.Range("L2").Value = 1594.5555
$ws.Range("M2").Value = -1684.2307
$ws.Range("N2").Value = -1820.5555
$ws.Range("H32").Value = 3079.3457
$ws.Range("I32").Value = 1671.9275
$ws.Range("K32").Value = 1671.9275
$ws.Range("M32").Value = -1384.9275
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31622
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -98112
$ws.Range("H74").Value = 1623.6415
$ws.Range("I74").Value = 1483.6216
$ws.Range("J74").Value = 1947.4375
$ws.Range("K74").Value = 1483.6216
$ws.Range("L74").Value = 1947.4375
$ws.Range("M74").Value = -609.6215999999999
$ws.Range("N74").Value = -3695.4375
$ws.Range("H77").Value = 1623.6415
$ws.Range("I77").Value = 1483.6216
$ws.Range("J77").Value = 1947.4375
$ws.Range("K77").Value = 7418.108
$ws.Range("L77").Value = 9737.1875
$ws.Range("M77").Value = -3050.108
$ws.Range("N77").Value = -18473.1875
$ws.Range("H97").Value = 38462576
$ws.Range("I97").Value = 45455476
$ws.Range("J97").Value = 1615.25
$ws.Range("K97").Value = 45455476
$ws.Range("L97").Value = 1615.25
$ws.Range("M97").Value = -45454980
$ws.Range("N97").Value = -2607.25
$ws.Range("H116").Value = 1745.1143
$ws.Range("I116").Value = 1797.2307
$ws.Range("J116").Value = 1594.5555
$ws.Range("K116").Value = 1797.2307
$ws.Range("L116").Value = 1594.5555
$ws.Range("M116").Value = 496.7692999999999
$ws.Range("N116").Value = -6182.5555
$ws.Range("H122").Value = 1687.7037
$ws.Range("I122").Value = 1539.4546
$ws.Range("J122").Value = 2340
$ws.Range("K122").Value = 4618.3638
$ws.Range("L122").Value = 7020
$ws.Range("M122").Value = -2168.3638
$ws.Range("N122").Value = -11920
$ws.Range("H133").Value = 40757.145
$ws.Range("J133").Value = 40757.145
$ws.Range("L133").Value = 40757.145
$ws.Range("N133").Value = -45817.145
$ws.Range("H135").Value = 24217.133
$ws.Range("J135").Value = 24217.133
$ws.Range("L135").Value = 24217.133
$ws.Range("N135").Value = -34357.133
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1745.1143
$ws.Range("I3").Value = 1797.2307
$ws.Range("J3").Value = 1594.5555
$ws.Range("K3").Value = 1797.2307
$ws.Range("L3").Value = 1594.5555
$ws.Range("M3").Value = -1683.2307
$ws.Range("N3").Value = -1822.5555
$ws.Range("H94").Value = 947
$ws.Range("I94").Value = 860.3226
$ws.Range("J94").Value = 1245.5555
$ws.Range("K94").Value = 860.3226
$ws.Range("L94").Value = 1245.5555
$ws.Range("M94").Value = -409.3226
$ws.Range("N94").Value = -2147.5555
$ws.Range("H99").Value = 2283.1333
$ws.Range("I99").Value = 2121.818
$ws.Range("J99").Value = 2726.75
$ws.Range("K99").Value = 2121.818
$ws.Range("L99").Value = 2726.75
$ws.Range("M99").Value = -623.8180000000002
$ws.Range("N99").Value = -5722.75
$ws.Range("H105").Value = 3530.0667
$ws.Range("I105").Value = 2402.5
$ws.Range("J105").Value = 3940.0908
$ws.Range("K105").Value = 2402.5
$ws.Range("L105").Value = 3940.0908
$ws.Range("M105").Value = -655.5
$ws.Range("N105").Value = -7434.0908
$ws.Range("H134").Value = 2723.1
$ws.Range("I134").Value = 1880.5172
$ws.Range("K134").Value = 5641.5516
$ws.Range("M134").Value = -3106.5516
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6807604
$ws.Range("I31").Value = 1580.6538
$ws.Range("K31").Value = 1580.6538
$ws.Range("M31").Value = -1285.6538
$ws.Range("H34").Value = 6807604
$ws.Range("I34").Value = 1580.6538
$ws.Range("K34").Value = 1580.6538
$ws.Range("M34").Value = -1378.6538
$ws.Range("H132").Value = 63909.22
$ws.Range("I132").Value = 1482.091
$ws.Range("J132").Value = 121134.086
$ws.Range("K132").Value = 4446.272999999999
$ws.Range("L132").Value = 363402.258
$ws.Range("M132").Value = -1916.272999999999
$ws.Range("N132").Value = -368462.258
$ws.Range("H134").Value = 431871.62
$ws.Range("I134").Value = 522954.4
$ws.Range("J134").Value = 158623.33
$ws.Range("K134").Value = 1568863.2
$ws.Range("L134").Value = 475869.99
$ws.Range("M134").Value = -1566328.2
$ws.Range("N134").Value = -480939.99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1445.3889
$ws.Range("I117").Value = 154.22223
$ws.Range("J117").Value = 2736.5557
$ws.Range("K117").Value = 462.66669
$ws.Range("L117").Value = 8209.667099999999
$ws.Range("M117").Value = 2979.33331
$ws.Range("N117").Value = -15093.6671
$ws.Range("H121").Value = 183124.64
$ws.Range("J121").Value = 282863.53
$ws.Range("L121").Value = 848590.5900000001
$ws.Range("N121").Value = -851210.5900000001
$ws.Range("H131").Value = 1093.25
$ws.Range("I131").Value = 486.66666
$ws.Range("J131").Value = 1131.9681
$ws.Range("K131").Value = 1459.99998
$ws.Range("L131").Value = 3395.9043
$ws.Range("M131").Value = 3580.00002
$ws.Range("N131").Value = -13475.9043
$ws.Range("H140").Value = 160021.84
$ws.Range("I140").Value = 215851.08
$ws.Range("K140").Value = 647553.24
$ws.Range("M140").Value = -642373.24
$ws.Range("H141").Value = 62505164
$ws.Range("I141").Value = 90914456
$ws.Range("K141").Value = 272743368
$ws.Range("M141").Value = -272738188
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3309.6086
$ws.Range("I97").Value = 2328.3333
$ws.Range("J97").Value = 6842.2
$ws.Range("K97").Value = 2328.3333
$ws.Range("L97").Value = 6842.2
$ws.Range("M97").Value = -1832.3333
$ws.Range("N97").Value = -7834.2
$ws.Range("H126").Value = 90920264
$ws.Range("I126").Value = 125014424
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 375043272
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -375040802
$ws.Range("N126").Value = -12440
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3641.125
$ws.Range("I40").Value = 3604.1428
$ws.Range("J40").Value = 3900
$ws.Range("K40").Value = 3604.1428
$ws.Range("L40").Value = 3900
$ws.Range("M40").Value = -3468.1428
$ws.Range("N40").Value = -4172
$ws.Range("H82").Value = 4630880
$ws.Range("I82").Value = 1144.4445
$ws.Range("J82").Value = 9260615
$ws.Range("K82").Value = 1144.4445
$ws.Range("L82").Value = 9260615
$ws.Range("M82").Value = -783.4445000000001
$ws.Range("N82").Value = -9261337
$ws.Range("H85").Value = 4630880
$ws.Range("I85").Value = 1144.4445
$ws.Range("J85").Value = 9260615
$ws.Range("K85").Value = 1144.4445
$ws.Range("L85").Value = 9260615
$ws.Range("M85").Value = 103.5554999999999
$ws.Range("N85").Value = -9263111
$ws.Range("H93").Value = 1156.4375
$ws.Range("I93").Value = 1249.8334
$ws.Range("J93").Value = 1100.4
$ws.Range("K93").Value = 1249.8334
$ws.Range("L93").Value = 1100.4
$ws.Range("M93").Value = -1.833399999999983
$ws.Range("N93").Value = -3596.4
$ws.Range("H100").Value = 2100
$ws.Range("I100").Value = 2150
$ws.Range("K100").Value = 2150
$ws.Range("M100").Value = -1609
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1350.1428
$ws.Range("I81").Value = 1119.6666
$ws.Range("J81").Value = 1523
$ws.Range("K81").Value = 2239.3332
$ws.Range("L81").Value = 3046
$ws.Range("M81").Value = -1178.3332
$ws.Range("N81").Value = -5168
$ws.Range("H84").Value = 1350.1428
$ws.Range("I84").Value = 1119.6666
$ws.Range("J84").Value = 1523
$ws.Range("K84").Value = 11196.666
$ws.Range("L84").Value = 15230
$ws.Range("M84").Value = -5892.666000000001
$ws.Range("N84").Value = -25838
$ws.Range("H122").Value = 1682310.4
$ws.Range("I122").Value = 1787392.2
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5362176.6
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -5359726.6
$ws.Range("N122").Value = -7900
